$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.807.67"
$ws.Range("E2").Value = "  +1.37%  "
$ws.Range("D3").Value = "1.700.04"
$ws.Range("E3").Value = "  +0.70%  "
$ws.Range("E4").Value = "  -0.37%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.43"
$ws.Range("E5").Value = "  +0.66%  "
$ws.Range("E6").Value = "  -0.21%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3989"
$ws.Range("E7").Value = "  +2.69%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4041"
$ws.Range("E8").Value = "  +0.21%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.001"
$ws.Range("E9").Value = "  -0.44%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.472"
$ws.Range("E10").Value = "  -1.42%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.64"
$ws.Range("E11").Value = "  +1.75%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08809"
$ws.Range("E12").Value = "  +0.71%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "26.08"
$ws.Range("E13").Value = "  +3.55%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.547"
$ws.Range("E14").Value = "  +0.73%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.001"
$ws.Range("E15").Value = "  +0.71%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001348"
$ws.Range("E16").Value = "  -0.08%  "
$ws.Range("D17").Value = "1.701.92"
$ws.Range("E17").Value = "  +0.95%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "95.78"
$ws.Range("E18").Value = "  -2.73%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07178"
$ws.Range("E19").Value = "  +1.26%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "20.95"
$ws.Range("E20").Value = "  +4.91%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.351"
$ws.Range("E21").Value = "  +1.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.003"
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.41"
$ws.Range("E23").Value = "  +1.32%  "
$ws.Range("D24").Value = "24.775.00"
$ws.Range("E24").Value = "  +1.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.375"
$ws.Range("E25").Value = "  +0.79%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.917"
$ws.Range("E26").Value = "  -1.88%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.20"
$ws.Range("E27").Value = "  +2.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.182"
$ws.Range("E28").Value = "  +18.46%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "161.71"
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "144.44"
$ws.Range("E30").Value = "  +5.38%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.329"
$ws.Range("E31").Value = "  -4.70%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.306"
$ws.Range("E32").Value = "  +17.20%  "
$ws.Range("D33").Value = "1.891.02"
$ws.Range("E33").Value = "  +1.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08646"
$ws.Range("E34").Value = "  -1.84%  "
$ws.Range("B35").Value = "InternetComputer(DFINITY)"
$ws.Range("C35").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.353"
$ws.Range("E35").Value = "  -0.33%  "
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.03185"
$ws.Range("E36").Value = "  +9.39%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.033"
$ws.Range("E37").Value = "  +0.23%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2839"
$ws.Range("E38").Value = "  +2.22%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.8340"
$ws.Range("E39").Value = "  +5.42%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "10.78"
$ws.Range("E40").Value = "  +0.13%  "
$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.09465"
$ws.Range("E41").Value = "  +3.65%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "14.19"
$ws.Range("E42").Value = "  -0.19%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.478"
$ws.Range("E43").Value = "  +1.41%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "17.75"
$ws.Range("E44").Value = "  +6.22%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.714"
$ws.Range("E45").Value = "  +4.31%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.7453"
$ws.Range("E46").Value = "  +3.33%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.224"
$ws.Range("E47").Value = "  +0.69%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.385"
$ws.Range("E48").Value = "  +2.59%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.001"
$ws.Range("E49").Value = "  -0.05%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.08380"
$ws.Range("E50").Value = "  +4.54%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "139.65"
$ws.Range("E51").Value = "  +1.29%  "
